# Update cryptocurrency price/volume figures per the Sep 30 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.949.64"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.673.42"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "1.910.29"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.689.93"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "26.966.78"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "1.479.79"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0172"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.50%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "1.817.08"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.775"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
